$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: note about the "servicios públicos" category being handled independently for now
$ws.Range("D9").Value = "* Por ahora esta categoría se maneja de manera independiente"

# Row 6: "localidades" becomes "ubicaciones" with its own tag NE00U00, plus a note in D
$ws.Range("A6").Value = "NE00U00"

# Row 5: "ciudades" now has its own dictionary tag NE00C00 (was incorrectly NE00P00)
$ws.Range("A5").Value = "NE00C00"

$ws.Range("B6").Value = "ubicaciones"
$ws.Range("D6").Value = "*No tenemos diccionarios de ubicaciones"

# Update the active selection to D7 (matches author's cursor position after edits)
$ws.Range("D7").Select()
